$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.059.75"
$ws.Range("E2").Value = "  +2.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.195.63"
$ws.Range("E3").Value = "  +1.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.60"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.37"
$ws.Range("E6").Value = "  +4.23%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("E10").Value = "  +2.51%  "

$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.746.29"
$ws.Range("E12").Value = "  +1.30%  "

$ws.Range("E13").Value = "  -1.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.86"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.091.24"
$ws.Range("E16").Value = "  +2.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.224.65"
$ws.Range("E17").Value = "  +2.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.27"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.33"
$ws.Range("E19").Value = "  +2.08%  "

$ws.Range("E20").Value = "  +0.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.72"
$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.64"
$ws.Range("E26").Value = "  +4.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("E28").Value = "  +1.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.48"
$ws.Range("E29").Value = "  +1.78%  "

$ws.Range("E30").Value = "  +0.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").Value = "  +2.24%  "

$ws.Range("E33").Value = "  +2.46%  "

$ws.Range("E34").Value = "  +4.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.67"
$ws.Range("E35").Value = "  -1.54%  "

$ws.Range("E36").Value = "  +1.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.841.10"
$ws.Range("E37").Value = "  +7.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.43"
$ws.Range("E38").Value = "  +5.05%  "

$ws.Range("E39").Value = "  +2.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0310"
$ws.Range("E40").Value = "  +8.64%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.25"
$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("E43").Value = "  +2.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.720"
$ws.Range("E44").Value = "  +1.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.104"
$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.237.97"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.75"
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("E50").Value = "  +4.97%  "

$ws.Range("E51").Value = "  +0.03%  "
